$wb = $excel.ActiveWorkbook

# "SheetInvalid" gets its A2 cell populated with the text "invalid",
# matching the format/value used to flag invalid rows on the other sheets.
$wsInvalid = $wb.Worksheets.Item("SheetInvalid")
$wsInvalid.Range("A2").Value = "invalid"

# "Sheet1" selection moves from the old full-range selection to C1.
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("C1").Select()

# Restore "SheetInvalid" as the active/selected sheet tab.
$wsInvalid.Activate()
